{"js": "// Replace the 25 multiplication equations in the table with the new set.\nconst replacements = [\n  [\"47\u00d782=3854\", \"43\u00d716=688\"],\n  [\"12\u00d770=840\", \"17\u00d724=408\"],\n  [\"70\u00d751=3570\", \"41\u00d752=2132\"],\n  [\"33\u00d726=858\", \"26\u00d768=1768\"],\n  [\"43\u00d794=4042\", \"83\u00d784=6972\"],\n  [\"19\u00d730=570\", \"89\u00d747=4183\"],\n  [\"56\u00d767=3752\", \"46\u00d718=828\"],\n  [\"83\u00d761=5063\", \"69\u00d798=6762\"],\n  [\"12\u00d720=240\", \"16\u00d719=304\"],\n  [\"20\u00d770=1400\", \"35\u00d769=2415\"],\n  [\"67\u00d742=2814\", \"82\u00d723=1886\"],\n  [\"45\u00d777=3465\", \"58\u00d789=5162\"],\n  [\"25\u00d780=2000\", \"70\u00d728=1960\"],\n  [\"76\u00d734=2584\", \"67\u00d755=3685\"],\n  [\"75\u00d751=3825\", \"27\u00d730=810\"],\n  [\"26\u00d793=2418\", \"94\u00d736=3384\"],\n  [\"45\u00d757=2565\", \"89\u00d711=979\"],\n  [\"58\u00d713=754\", \"88\u00d773=6424\"],\n  [\"84\u00d718=1512\", \"93\u00d739=3627\"],\n  [\"63\u00d728=1764\", \"73\u00d723=1679\"],\n  [\"50\u00d756=2800\", \"72\u00d724=1728\"],\n  [\"43\u00d734=1462\", \"81\u00d779=6399\"],\n  [\"71\u00d748=3408\", \"46\u00d717=782\"],\n  [\"35\u00d776=2660\", \"34\u00d728=952\"],\n  [\"72\u00d783=5976\", \"56\u00d795=5320\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 multiplication equations in the table with the new set.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{Old=\"47\u00d782=3854\"; New=\"43\u00d716=688\"},\n    @{Old=\"12\u00d770=840\"; New=\"17\u00d724=408\"},\n    @{Old=\"70\u00d751=3570\"; New=\"41\u00d752=2132\"},\n    @{Old=\"33\u00d726=858\"; New=\"26\u00d768=1768\"},\n    @{Old=\"43\u00d794=4042\"; New=\"83\u00d784=6972\"},\n    @{Old=\"19\u00d730=570\"; New=\"89\u00d747=4183\"},\n    @{Old=\"56\u00d767=3752\"; New=\"46\u00d718=828\"},\n    @{Old=\"83\u00d761=5063\"; New=\"69\u00d798=6762\"},\n    @{Old=\"12\u00d720=240\"; New=\"16\u00d719=304\"},\n    @{Old=\"20\u00d770=1400\"; New=\"35\u00d769=2415\"},\n    @{Old=\"67\u00d742=2814\"; New=\"82\u00d723=1886\"},\n    @{Old=\"45\u00d777=3465\"; New=\"58\u00d789=5162\"},\n    @{Old=\"25\u00d780=2000\"; New=\"70\u00d728=1960\"},\n    @{Old=\"76\u00d734=2584\"; New=\"67\u00d755=3685\"},\n    @{Old=\"75\u00d751=3825\"; New=\"27\u00d730=810\"},\n    @{Old=\"26\u00d793=2418\"; New=\"94\u00d736=3384\"},\n    @{Old=\"45\u00d757=2565\"; New=\"89\u00d711=979\"},\n    @{Old=\"58\u00d713=754\"; New=\"88\u00d773=6424\"},\n    @{Old=\"84\u00d718=1512\"; New=\"93\u00d739=3627\"},\n    @{Old=\"63\u00d728=1764\"; New=\"73\u00d723=1679\"},\n    @{Old=\"50\u00d756=2800\"; New=\"72\u00d724=1728\"},\n    @{Old=\"43\u00d734=1462\"; New=\"81\u00d779=6399\"},\n    @{Old=\"71\u00d748=3408\"; New=\"46\u00d717=782\"},\n    @{Old=\"35\u00d776=2660\"; New=\"34\u00d728=952\"},\n    @{Old=\"72\u00d783=5976\"; New=\"56\u00d795=5320\"},\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute($pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null\n}\n"}
